# Add a new worksheet ("Sheet1") after the last existing sheet, make it the
# active sheet, and populate its header row with the "Panel" cable/connector
# table columns: Cable, Component, AWG, Pin Count, Length, Conn A, Conn B, Wrap.

$wb = $excel.ActiveWorkbook

# The previously-active sheet ("Upper Level Wiring") had its cursor moved
# before the new sheet was created; update its remembered selection too.
$prevActive = $wb.Worksheets.Item(2)
[void]$prevActive.Range("I20").Select()

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)

# Make it the active / selected tab (also clears tabSelected on the old sheet).
$ws.Activate()

# Column widths to roughly match the authored layout (col B and col C).
$ws.Columns.Item(2).ColumnWidth = 11.666666666666668
$ws.Columns.Item(3).ColumnWidth = 6.333333333333334

# Fill header row right-to-left so new shared-string entries are appended in
# the same order as the source workbook (H1, G1, F1, E1, D1, then A1; B1/C1
# reuse strings that already exist in the shared string table).
$ws.Range("H1").Value = "Wrap"
$ws.Range("G1").Value = "Conn B"
$ws.Range("F1").Value = "Conn A"
$ws.Range("E1").Value = "Length"
$ws.Range("D1").Value = "Pin Count"
$ws.Range("C1").Value = "AWG"
$ws.Range("B1").Value = "Component"
$ws.Range("A1").Value = "Cable"

# Match the selection left behind in the source file.
[void]$ws.Range("I6").Select()
